# Apply the latest cryptos snapshot (GitHub Actions scheduled refresh) to Sheet1.
# Columns: B=Coin name, C=Coinranking link, D=Price, E=Volume(1h) change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Coin = "Bitcoin"; Link = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; Price = "24.636.64"; Volume = "  -1.46%  " },
    @{ Row = 3; Coin = "Ethereum"; Link = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; Price = "1.674.65"; Volume = "  -2.13%  " },
    @{ Row = 4; Coin = "TetherUSD"; Link = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; Price = "1.004"; Volume = "  +0.54%  " },
    @{ Row = 5; Coin = "BNB"; Link = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; Price = "314.32"; Volume = "  -1.11%  " },
    @{ Row = 6; Coin = "USDC"; Link = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; Price = "1.004"; Volume = "  +0.46%  " },
    @{ Row = 7; Coin = "XRP"; Link = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; Price = "0.3908"; Volume = "  -3.43%  " },
    @{ Row = 8; Coin = "Cardano"; Link = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; Price = "0.3937"; Volume = "  -3.66%  " },
    @{ Row = 9; Coin = "OKB"; Link = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; Price = "52.02"; Volume = "  -3.74%  " },
    @{ Row = 10; Coin = "BinanceUSD"; Link = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; Price = "1.005"; Volume = "  +0.25%  " },
    @{ Row = 11; Coin = "Polygon"; Link = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; Price = "1.388"; Volume = "  -6.40%  " },
    @{ Row = 12; Coin = "Dogecoin"; Link = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; Price = "0.08635"; Volume = "  -2.45%  " },
    @{ Row = 13; Coin = "Solana"; Link = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; Price = "25.13"; Volume = "  -4.98%  " },
    @{ Row = 14; Coin = "Polkadot"; Link = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; Price = "7.304"; Volume = "  -3.06%  " },
    @{ Row = 15; Coin = "Chainlink"; Link = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; Price = "7.740"; Volume = "  -5.01%  " },
    @{ Row = 16; Coin = "ShibaInu"; Link = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; Price = "0.00001313"; Volume = "  -3.55%  " },
    @{ Row = 17; Coin = "WrappedEther"; Link = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; Price = "1.675.45"; Volume = "  -4.52%  " },
    @{ Row = 18; Coin = "Litecoin"; Link = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; Price = "93.64"; Volume = "  -3.66%  " },
    @{ Row = 19; Coin = "TRON"; Link = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; Price = "0.07064"; Volume = "  -1.43%  " },
    @{ Row = 20; Coin = "Avalanche"; Link = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; Price = "20.44"; Volume = "  -3.82%  " },
    @{ Row = 21; Coin = "Uniswap"; Link = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; Price = "7.032"; Volume = "  -3.61%  " },
    @{ Row = 22; Coin = "Dai"; Link = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; Price = "1.004"; Volume = "  +0.55%  " },
    @{ Row = 23; Coin = "Cosmos"; Link = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Price = "13.90"; Volume = "  -3.69%  " },
    @{ Row = 24; Coin = "WrappedBTC"; Link = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; Price = "24.646.35"; Volume = "  -1.37%  " },
    @{ Row = 25; Coin = "Toncoin"; Link = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; Price = "2.352"; Volume = "  +1.20%  " },
    @{ Row = 26; Coin = "EthereumClassic"; Link = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; Price = "23.18"; Volume = "  -0.67%  " },
    @{ Row = 27; Coin = "LidoDAOToken"; Link = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; Price = "2.738"; Volume = "  -6.65%  " },
    @{ Row = 28; Coin = "Monero"; Link = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; Price = "162.32"; Volume = "  -2.88%  " },
    @{ Row = 29; Coin = "HuobiToken"; Link = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; Price = "5.734"; Volume = "  -7.73%  " },
    @{ Row = 30; Coin = "BitcoinCash"; Link = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; Price = "146.56"; Volume = "  -0.41%  " },
    @{ Row = 31; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "7.862"; Volume = "  -6.91%  " },
    @{ Row = 32; Coin = "WEMIXTOKEN"; Link = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; Price = "2.448"; Volume = "  +9.01%  " },
    @{ Row = 33; Coin = "WrappedliquidstakedEther2.0"; Link = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; Price = "1.853.19"; Volume = "  -0.55%  " },
    @{ Row = 34; Coin = "Hedera"; Link = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; Price = "0.08380"; Volume = "  -5.85%  " },
    @{ Row = 35; Coin = "VeChain"; Link = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; Price = "0.03030"; Volume = "  -6.10%  " },
    @{ Row = 36; Coin = "InternetComputer(DFINITY)"; Link = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; Price = "6.898"; Volume = "  -5.35%  " },
    @{ Row = 37; Coin = "Algorand"; Link = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; Price = "0.2794"; Volume = "  -2.56%  " },
    @{ Row = 38; Coin = "ImmutableX"; Link = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; Price = "0.9850"; Volume = "  -4.51%  " },
    @{ Row = 39; Coin = "Stellar"; Link = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; Price = "0.09464"; Volume = "  +0.79%  " },
    @{ Row = 40; Coin = "FraxShare"; Link = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; Price = "10.54"; Volume = "  -3.39%  " },
    @{ Row = 41; Coin = "TrustWalletToken"; Link = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; Price = "1.531"; Volume = "  +4.06%  " },
    @{ Row = 42; Coin = "TheSandbox"; Link = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; Price = "0.7864"; Volume = "  -7.50%  " },
    @{ Row = 43; Coin = "Aptos"; Link = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; Price = "13.49"; Volume = "  -5.25%  " },
    @{ Row = 44; Coin = "EnergySwap"; Link = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Price = "16.37"; Volume = "  -6.21%  " },
    @{ Row = 45; Coin = "Decentraland"; Link = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; Price = "0.7113"; Volume = "  -4.60%  " },
    @{ Row = 46; Coin = "NEARProtocol"; Link = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; Price = "2.551"; Volume = "  -6.22%  " },
    @{ Row = 47; Coin = "PancakeSwap"; Link = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; Price = "4.189"; Volume = "  -1.51%  " },
    @{ Row = 48; Coin = "Cronos"; Link = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; Price = "0.08639"; Volume = "  +3.07%  " },
    @{ Row = 49; Coin = "Frax"; Link = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"; Price = "1.003"; Volume = "  +0.47%  " },
    @{ Row = 50; Coin = "Flow"; Link = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"; Price = "1.323"; Volume = "  -5.67%  " },
    @{ Row = 51; Coin = "Quant"; Link = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; Price = "137.31"; Volume = "  -3.52%  " }
)

foreach ($r in $rows) {
    $ws.Range("B$($r.Row)").Value = $r.Coin
    $ws.Range("C$($r.Row)").Value = $r.Link

    # Price strings look numeric (e.g. "1.004") and Excel would silently
    # coerce them to a Double on assignment, dropping the trailing zero /
    # the multi-dot thousands grouping Coinranking uses (e.g. "24.636.64").
    # Force the cell to Text first so the literal string is stored verbatim,
    # then drop the explicit format again so no stray style is left behind.
    $priceCell = $ws.Range("D$($r.Row)")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $r.Price
    $priceCell.ClearFormats()

    $ws.Range("E$($r.Row)").Value = $r.Volume
}
